$wb = $excel.ActiveWorkbook

# This script applies updated market-price / profit figures produced by
# the scheduled pricing runner. Each worksheet (one per crafting class)
# has static columns H:N holding price/profit snapshots for specific rows;
# these are refreshed in place, cell by cell, per the runner output.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I5").Value = 618.7778
$ws.Range("K5").Value = 618.7778
$ws.Range("H5").Value = 561.9
$ws.Range("M5").Value = -503.7778
$ws.Range("I40").Value = 4973.25
$ws.Range("K40").Value = 4973.25
$ws.Range("H40").Value = 5565385
$ws.Range("M40").Value = -4798.25
$ws.Range("I70").Value = 3449.8572
$ws.Range("K70").Value = 10349.5716
$ws.Range("H70").Value = 7294.7
$ws.Range("M70").Value = -10079.5716
$ws.Range("I73").Value = 3449.8572
$ws.Range("K73").Value = 10349.5716
$ws.Range("H73").Value = 7294.7
$ws.Range("M73").Value = -9413.571599999999
$ws.Range("I98").Value = 1604.6111
$ws.Range("K98").Value = 1604.6111
$ws.Range("H98").Value = 1604.6111
$ws.Range("M98").Value = -106.6111000000001
$ws.Range("I122").Value = 1604.6111
$ws.Range("K122").Value = 4813.8333
$ws.Range("H122").Value = 1604.6111
$ws.Range("M122").Value = -2363.8333
$ws.Range("I129").Value = 1539.7142
$ws.Range("K129").Value = 4619.142599999999
$ws.Range("N129").Value = -17399.5
$ws.Range("H129").Value = 1876.7273
$ws.Range("J129").Value = 2466.5
$ws.Range("L129").Value = 7399.5
$ws.Range("M129").Value = 380.8574000000008
$ws.Range("I135").Value = 38462744
$ws.Range("K135").Value = 346164696
$ws.Range("H135").Value = 107143976
$ws.Range("M135").Value = -346162161

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I61").Value = 100006860
$ws.Range("K61").Value = 100006860
$ws.Range("H61").Value = 83339720
$ws.Range("M61").Value = -100006648
$ws.Range("N74").Value = -3623
$ws.Range("H74").Value = 55562372
$ws.Range("J74").Value = 1875
$ws.Range("L74").Value = 1875
$ws.Range("N77").Value = -18111
$ws.Range("H77").Value = 55562372
$ws.Range("J77").Value = 1875
$ws.Range("L77").Value = 9375
$ws.Range("I136").Value = 100006860
$ws.Range("K136").Value = 300020580
$ws.Range("H136").Value = 83339720
$ws.Range("M136").Value = -300018030

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 5943.05
$ws.Range("K31").Value = 5943.05
$ws.Range("H31").Value = 10019.789
$ws.Range("M31").Value = -5648.05
$ws.Range("I34").Value = 5943.05
$ws.Range("K34").Value = 5943.05
$ws.Range("H34").Value = 10019.789
$ws.Range("M34").Value = -5741.05
$ws.Range("I41").Value = 11273.667
$ws.Range("K41").Value = 11273.667
$ws.Range("H41").Value = 18685.572
$ws.Range("M41").Value = -10845.667
$ws.Range("I47").Value = 17942.5
$ws.Range("K47").Value = 17942.5
$ws.Range("H47").Value = 25294.666
$ws.Range("M47").Value = -17376.5
$ws.Range("N50").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N51").Value = -60972
$ws.Range("H51").Value = 59500
$ws.Range("J51").Value = 59500
$ws.Range("L51").Value = 59500
$ws.Range("I58").Value = 62513296
$ws.Range("K58").Value = 62513296
$ws.Range("N58").Value = -6739
$ws.Range("H58").Value = 45465940
$ws.Range("J58").Value = 6333
$ws.Range("L58").Value = 6333
$ws.Range("M58").Value = -62513093
$ws.Range("N59").Value = -133956.67
$ws.Range("H59").Value = 123776
$ws.Range("J59").Value = 131666.67
$ws.Range("L59").Value = 131666.67
$ws.Range("N60").Value = -25041.8
$ws.Range("H60").Value = 22661.842
$ws.Range("J60").Value = 24019.8
$ws.Range("L60").Value = 24019.8
$ws.Range("N61").Value = -60196
$ws.Range("H61").Value = 59500
$ws.Range("J61").Value = 59500
$ws.Range("L61").Value = 59500
$ws.Range("N68").Value = -61498
$ws.Range("H68").Value = 80000
$ws.Range("J68").Value = 60000
$ws.Range("L68").Value = 60000
$ws.Range("N71").Value = -187488
$ws.Range("H71").Value = 80000
$ws.Range("J71").Value = 60000
$ws.Range("L71").Value = 180000
$ws.Range("I105").Value = 2858259.2
$ws.Range("K105").Value = 2858259.2
$ws.Range("H105").Value = 1819664.9
$ws.Range("M105").Value = -2856512.2
$ws.Range("N107").Value = -72947.53
$ws.Range("H107").Value = 595419.6
$ws.Range("J107").Value = 69107.53
$ws.Range("L107").Value = 69107.53
$ws.Range("I134").Value = 125001750
$ws.Range("K134").Value = 375005250
$ws.Range("H134").Value = 83334870
$ws.Range("M134").Value = -375002715
$ws.Range("I136").Value = 62513296
$ws.Range("K136").Value = 187539888
$ws.Range("N136").Value = -24099
$ws.Range("H136").Value = 45465940
$ws.Range("J136").Value = 6333
$ws.Range("L136").Value = 18999
$ws.Range("M136").Value = -187537338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I8").Value = 208.5
$ws.Range("K8").Value = 625.5
$ws.Range("H8").Value = 208.5
$ws.Range("M8").Value = -486.5
$ws.Range("I86").Value = 176
$ws.Range("K86").Value = 528
$ws.Range("N86").Value = -6844.5002
$ws.Range("H86").Value = 964.9
$ws.Range("J86").Value = 1490.8334
$ws.Range("L86").Value = 4472.5002
$ws.Range("M86").Value = 658
$ws.Range("I89").Value = 176
$ws.Range("K89").Value = 1584
$ws.Range("N89").Value = -25273.5006
$ws.Range("H89").Value = 964.9
$ws.Range("J89").Value = 1490.8334
$ws.Range("L89").Value = 13417.5006
$ws.Range("M89").Value = 4344
$ws.Range("I121").Value = 339999.66
$ws.Range("K121").Value = 1019998.98
$ws.Range("N121").Value = -60163.819
$ws.Range("H121").Value = 87928.07000000001
$ws.Range("J121").Value = 19181.273
$ws.Range("L121").Value = 57543.819
$ws.Range("M121").Value = -1018688.98

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N53").Value = -51261
$ws.Range("H53").Value = 49999
$ws.Range("J53").Value = 49999
$ws.Range("L53").Value = 49999
$ws.Range("I122").Value = 68706.78
$ws.Range("K122").Value = 206120.34
$ws.Range("H122").Value = 57589.184
$ws.Range("M122").Value = -203670.34
$ws.Range("I126").Value = 4141.4165
$ws.Range("K126").Value = 12424.2495
$ws.Range("H126").Value = 4164.6216
$ws.Range("M126").Value = -9954.249500000002
$ws.Range("I132").Value = 6251360.5
$ws.Range("K132").Value = 18754081.5
$ws.Range("H132").Value = 5209797.5
$ws.Range("M132").Value = -18751551.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I58").Value = 291984.28
$ws.Range("K58").Value = 291984.28
$ws.Range("H58").Value = 232742.11
$ws.Range("M58").Value = -291724.28

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 5501
$ws.Range("K62").Value = 5501
$ws.Range("N62").Value = -20578
$ws.Range("H62").Value = 13798.4
$ws.Range("J62").Value = 19330
$ws.Range("L62").Value = 19330
$ws.Range("M62").Value = -4877
$ws.Range("I65").Value = 5501
$ws.Range("K65").Value = 27505
$ws.Range("N65").Value = -102890
$ws.Range("H65").Value = 13798.4
$ws.Range("J65").Value = 19330
$ws.Range("L65").Value = 96650
$ws.Range("M65").Value = -24385
$ws.Range("I136").Value = 45455424
$ws.Range("K136").Value = 136366272
$ws.Range("H136").Value = 41667972
$ws.Range("M136").Value = -136363722
$ws.Range("N138").Value = -107992.43
$ws.Range("H138").Value = 97712.42999999999
$ws.Range("J138").Value = 97712.42999999999
$ws.Range("L138").Value = 97712.42999999999
